$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.264.80'
$ws.Range("D3").Value = '1.594.49'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +0.24%  '
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.98'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0855'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("D12").Value = '1.818.90'
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("D13").Value = '1.585.16'
$ws.Range("E13").Value = '  -0.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.98'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.47'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("D17").Value = '26.258.37'
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.66'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.03%  '
$ws.Range("E20").Value = '  -0.68%  '
$ws.Range("E24").Value = '  -0.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.18'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.25%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("E29").Value = '  +1.73%  '
$ws.Range("E30").Value = '  -0.25%  '
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("B32").Value = 'Maker'
$ws.Range("C32").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D32").Value = '1.473.86'
$ws.Range("E32").Value = '  +4.42%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.19'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.53%  '
$ws.Range("E34").Value = '  -0.58%  '
$ws.Range("E35").Value = '  -0.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.566'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.50%  '
$ws.Range("E38").Value = '  -0.70%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.818'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.75'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.36%  '
$ws.Range("E42").Value = '  +1.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.932'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.31%  '
$ws.Range("D44").Value = '1.732.29'
$ws.Range("E44").Value = '  +0.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.755'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.36'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.83%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.25'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.71%  '
$ws.Range("E48").Value = '  -0.95%  '
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0951'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.38%  '
$ws.Range("E51").Value = '  -0.09%  '
